# Update calibration data with new costs.
# Rows 100-107 and 114-115 on the active sheet have columns J:AS filled with a
# single repeated value per row; replace that repeated value with the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    100 = 98191.44339
    101 = 319413.473
    102 = 174443.0271
    103 = 50278.36372
    104 = 75531.87953000001
    105 = 24248.29007
    106 = 19188.73298
    107 = 657263.9674
    114 = 1270.335465
    115 = 1207365.632
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Range("J$row`:AS$row").Value = $value
}
